# Threshold/Zn/2His_1Asp/5.xlsx — update geometric-parameter threshold table.
#
# The "theta_threshold_range" row is removed entirely (its row, its shared
# string, and the now-unused shared-string slot all go away), which shifts
# "pie_threshold_range" up by one row. The alpha/beta distance ranges and
# the (new) pie range also get new Min/Max numbers, and the view is left
# with the selection on C4 and column C very slightly narrower.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "theta_threshold_range" row (row 5). Deleting the whole row
# shifts "pie_threshold_range" up from row 6 to row 5 and drops the
# now-unused "theta_threshold_range" shared string automatically.
$ws.Rows(5).Delete()

# alpha_distance_range (row 2): Min/Max
$ws.Cells.Item(2, 2).Value = 3.8
$ws.Cells.Item(2, 3).Value = 11.9

# beta_distance_range (row 3): Min/Max
$ws.Cells.Item(3, 2).Value = 3.5
$ws.Cells.Item(3, 3).Value = 10.1

# ratio_threshold_range (row 4): unchanged (0.8 / 1.4)

# pie_threshold_range (now row 5 after the delete above): Min/Max
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = 20

# Column C is very slightly narrower than before (bestFit 5.875 -> ~5.5 chars)
$ws.Columns(3).ColumnWidth = 4.8

# Saved selection moves to C4
$ws.Range("C4").Select() | Out-Null

# Page setup: paper size A4, portrait orientation
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
